$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 46 ("Ministry of Economy and Trade"),
# shifting rows 46-52 down by one (46->47, 47->48, 48->49, 51->52, 52->53).
$ws.Rows("46:46").Insert()

# Remove the hyperlink that was previously anchored on the old A47
# (now shifted to A48) before we overwrite its text/style.
$ws.Hyperlinks.Delete()

# Row 46 (new, blank) becomes an empty "source" line.
$ws.Range("A46").Value = ""
$ws.Range("A46").Style = "source"

# Row 47 keeps "Ministry of Economy and Trade" but loses the hyperlink style,
# becoming plain "source" style text instead.
$ws.Range("A47").Value = "Ministry of Economy and Trade"
$ws.Range("A47").Style = "source"

# Row 48 becomes the blank "source" line (it used to hold the URL).
$ws.Range("A48").Value = ""
$ws.Range("A48").Style = "source"

# Row 49 now holds the URL text (plain, no hyperlink), "source" style.
$ws.Range("A49").Value = "http://www.economy.gov.lb/public/uploads/files/9524_6086_6462.pdf"
$ws.Range("A49").Style = "source"

# Row 52 unchanged: "CAS - MET" in "title" style.
$ws.Range("A52").Value = "CAS - MET"
$ws.Range("A52").Style = "title"

# Row 53 (previously held the long citation text) now just repeats
# "CAS - MET" in "source" style.
$ws.Range("A53").Value = "CAS - MET"
$ws.Range("A53").Style = "source"
